$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the bot text in cell A39: "тексты" -> "туры"
$ws.Range("A39").Value = "Вот туры бота, исправте их и отошлите их мне"

# Reflect the user's scroll/selection position at time of edit
$ws.Range("A39").Select()
$excel.ActiveWindow.ScrollRow = 32
